$d = $word.ActiveDocument

$replacements = @(
    @{ old = "348÷6="; new = "724÷6=" },
    @{ old = "409÷6="; new = "116÷9=" },
    @{ old = "985÷9="; new = "754÷7=" },
    @{ old = "918÷7="; new = "741÷2=" },
    @{ old = "867÷4="; new = "891÷6=" },
    @{ old = "345÷6="; new = "834÷8=" },
    @{ old = "946÷5="; new = "769÷3=" },
    @{ old = "944÷4="; new = "748÷7=" },
    @{ old = "411÷2="; new = "299÷2=" },
    @{ old = "546÷7="; new = "928÷9=" },
    @{ old = "154÷4="; new = "461÷6=" },
    @{ old = "807÷2="; new = "336÷4=" },
    @{ old = "535÷7="; new = "929÷8=" },
    @{ old = "959÷4="; new = "244÷5=" },
    @{ old = "261÷4="; new = "973÷7=" },
    @{ old = "653÷5="; new = "775÷4=" },
    @{ old = "470÷4="; new = "719÷3=" },
    @{ old = "486÷3="; new = "159÷6=" },
    @{ old = "330÷5="; new = "309÷9=" },
    @{ old = "319÷4="; new = "843÷8=" },
    @{ old = "532÷9="; new = "361÷4=" },
    @{ old = "269÷9="; new = "953÷8=" },
    @{ old = "444÷4="; new = "763÷9=" },
    @{ old = "389÷5="; new = "135÷2=" },
    @{ old = "899÷7="; new = "549÷2=" }
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
